$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7383365035057068
$ws.Range("B1").Value = 2.015712022781372
$ws.Range("C1").Value = 3.948572635650635
$ws.Range("D1").Value = 3.515429019927979
$ws.Range("E1").Value = 2.028938293457031
